# Add a new character style "Verbatim Char" (styleId "VerbatimChar"),
# based on the existing "Body Text Char" style, with its font set to
# Consolas — mirrors the reference diff's new <w:style> block appended
# right after BodyTextChar in styles.xml.

$d = $word.ActiveDocument

# wdStyleTypeCharacter = 2
$verbatimChar = $d.Styles.Add("VerbatimChar", 2)
$verbatimChar.NameLocal = "Verbatim Char"
$verbatimChar.BaseStyle = $d.Styles("BodyTextChar")
$verbatimChar.Font.Name = "Consolas"
